$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the erroneous "Paper # 4" rows (rows 4-6), which were all zeros
# and duplicated the axis-label bug described in the commit message.
# Deleting these rows shifts all subsequent rows up by 3.
$ws.Range("A4:A6").EntireRow.Delete()
